$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 and C2 values
$ws.Range("A2").Value = 7630
$ws.Range("C2").Value = 4

# Delete rows 3 and 4 entirely (content no longer present)
$ws.Range("A3:C4").Delete()
